$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.741.00'
$ws.Range("E2").Value = '  +1.45%  '

$ws.Range("D3").Value = '1.863.37'
$ws.Range("E3").Value = '  +1.37%  '

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.035'
$c.ClearFormats()
$ws.Range("E4").Value = '  +0.97%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '323.56'
$c.ClearFormats()
$ws.Range("E5").Value = '  +1.24%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '1.032'
$c.ClearFormats()
$ws.Range("E6").Value = '  +1.03%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.4423'
$c.ClearFormats()
$ws.Range("E7").Value = '  +1.77%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3798'
$c.ClearFormats()
$ws.Range("E8").Value = '  +2.14%  '

$ws.Range("E9").Value = '  +1.78%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.8863'
$c.ClearFormats()
$ws.Range("E10").Value = '  +1.29%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '21.76'
$c.ClearFormats()
$ws.Range("E11").Value = '  +1.95%  '

$ws.Range("D12").Value = '1.868.51'
$ws.Range("E12").Value = '  -8.39%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '5.549'
$c.ClearFormats()
$ws.Range("E13").Value = '  +1.48%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '6.769'
$c.ClearFormats()
$ws.Range("E14").Value = '  +1.41%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.07209'
$c.ClearFormats()
$ws.Range("E15").Value = '  +0.77%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '84.29'
$c.ClearFormats()
$ws.Range("E16").Value = '  +2.78%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '1.038'
$c.ClearFormats()
$ws.Range("E17").Value = '  +1.11%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.000009108'
$c.ClearFormats()
$ws.Range("E18").Value = '  +1.17%  '

$ws.Range("E19").Value = '  +1.19%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '15.57'
$c.ClearFormats()
$ws.Range("E20").Value = '  +1.06%  '

$ws.Range("D21").Value = '27.753.86'
$ws.Range("E21").Value = '  +1.35%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '5.306'
$c.ClearFormats()
$ws.Range("E22").Value = '  +1.30%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '11.33'
$c.ClearFormats()
$ws.Range("E23").Value = '  +2.13%  '

$ws.Range("D24").Value = '2.092.85'
$ws.Range("E24").Value = '  -5.98%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.015'
$c.ClearFormats()
$ws.Range("E25").Value = '  +5.75%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '158.80'
$c.ClearFormats()
$ws.Range("E26").Value = '  +1.44%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '18.84'
$c.ClearFormats()
$ws.Range("E27").Value = '  +1.83%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '1.996'
$c.ClearFormats()
$ws.Range("E28").Value = '  +3.44%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '5.336'
$c.ClearFormats()
$ws.Range("E29").Value = '  +1.06%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '118.06'
$c.ClearFormats()
$ws.Range("E30").Value = '  +2.33%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.09031'
$c.ClearFormats()
$ws.Range("E31").Value = '  +0.36%  '

$ws.Range("B32").Value = 'ARBITRUM'
$ws.Range("C32").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '1.222'
$c.ClearFormats()
$ws.Range("E32").Value = '  +2.11%  '

$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.7789'
$c.ClearFormats()
$ws.Range("E33").Value = '  +2.73%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '3.025'
$c.ClearFormats()
$ws.Range("E34").Value = '  +6.69%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '4.575'
$c.ClearFormats()
$ws.Range("E35").Value = '  +2.66%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '1.034'
$c.ClearFormats()
$ws.Range("E36").Value = '  +1.00%  '

$ws.Range("E37").Value = '  +0.03%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.01994'
$c.ClearFormats()
$ws.Range("E38").Value = '  +2.24%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.05348'
$c.ClearFormats()
$ws.Range("E39").Value = '  +1.74%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '2.873'
$c.ClearFormats()
$ws.Range("E40").Value = '  +2.69%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.5205'
$c.ClearFormats()
$ws.Range("E41").Value = '  +0.86%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.1693'
$c.ClearFormats()
$ws.Range("E42").Value = '  +1.92%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '6.877'
$c.ClearFormats()
$ws.Range("E43").Value = '  +5.58%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '8.678'
$c.ClearFormats()
$ws.Range("E44").Value = '  +2.70%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '110.53'
$c.ClearFormats()
$ws.Range("E45").Value = '  +2.47%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.06728'
$c.ClearFormats()
$ws.Range("E46").Value = '  +7.14%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '10.64'
$c.ClearFormats()
$ws.Range("E47").Value = '  +1.06%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '1.714'
$c.ClearFormats()
$ws.Range("E48").Value = '  +2.95%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.4718'
$c.ClearFormats()
$ws.Range("E49").Value = '  +2.11%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '1.915'
$c.ClearFormats()
$ws.Range("E50").Value = '  +1.02%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '39.77'
$c.ClearFormats()
$ws.Range("E51").Value = '  +1.64%  '

